# Bump the "Förändrad" (changed) date in column C for every data row
# (row 2 through the last used row) from 2023-09-23 (45192) to
# 2023-10-03 (45202).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45192) {
        $cell.Value2 = 45202
    }
}
